# [Outlook] (sensitivity label) Republish sensitivity label snippets (#768)
#
# Inserts 4 new "sensitivity label" API snippet rows at the end of the
# "Snippets" table (pushing the existing SessionData/Time rows down by 4
# rows), and grows the table / autofilter / used range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow the table (and its autofilter) so the new rows become part of it.
$lo.Resize($ws.Range("A1:F287"))

# Carry the existing row formatting (e.g. the numeric style on column D)
# down onto the brand-new rows 284-287 before we populate them.
$ws.Range("A283:F283").Copy()
$ws.Range("A284:F287").PasteSpecial(-4122)  # xlPasteFormats

function Set-Row([int]$r, [string]$a, [string]$b, [string]$c, $d, [string]$e, [string]$f) {
    $ws.Range("A$r").Value = $a
    $ws.Range("B$r").Value = $b
    $ws.Range("C$r").Value = $c
    $ws.Range("D$r").Value = $d
    $ws.Range("E$r").Value = $e
    $ws.Range("F$r").Value = $f
}

# New rows 276-279: SensitivityLabel / SensitivityLabelsCatalog snippets.
Set-Row 276 "Office" "SensitivityLabel" "getAsync" 2 "outlook-sensitivity-labels-sensitivity-label" "getCurrentSensitivityLabel"
Set-Row 277 "Office" "SensitivityLabel" "setAsync" 2 "outlook-sensitivity-labels-sensitivity-label" "setSensitivityLabel"
Set-Row 278 "Office" "SensitivityLabelsCatalog" "getAsync" 2 "outlook-sensitivity-labels-sensitivity-labels-catalog" "getSensitivityLabelsCatalog"
Set-Row 279 "Office" "SensitivityLabelsCatalog" "getIsEnabledAsync" 2 "outlook-sensitivity-labels-sensitivity-labels-catalog" "getSensitivityLabelsCatalogIsEnabled"

# Existing SessionData rows, shifted down by 4 (previously rows 276-280).
Set-Row 280 "Office" "SessionData" "clearAsync" 1 "outlook-event-based-activation-session-data-apis" "clearSessionData"
Set-Row 281 "Office" "SessionData" "getAllAsync" 1 "outlook-event-based-activation-session-data-apis" "getAllSessionData"
Set-Row 282 "Office" "SessionData" "getAsync" 1 "outlook-event-based-activation-session-data-apis" "getSessionData"
Set-Row 283 "Office" "SessionData" "removeAsync" 1 "outlook-event-based-activation-session-data-apis" "removeSessionData"
Set-Row 284 "Office" "SessionData" "setAsync" 1 "outlook-event-based-activation-session-data-apis" "setSessionData"

# Existing Time rows, shifted down by 4 (previously rows 281-283).
Set-Row 285 "Office" "Time" "getAsync" 2 "outlook-other-item-apis-get-set-start-appointment-organizer" "get"
Set-Row 286 "Office" "Time" "setAsync" 1 "outlook-other-item-apis-get-set-start-appointment-organizer" "set"
Set-Row 287 "Office" "Time" "setAsync" 1 "outlook-other-item-apis-get-set-end-appointment-organizer" "set"

# Match the updated selection left behind in the source file.
$ws.Range("F277").Select() | Out-Null
